$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value = -45.10157721253326
$ws.Range("C3").Value = -40.10157721253326
$ws.Range("D3").Value = -35.10157721253326
$ws.Range("E3").Value = -25.10157721253326
$ws.Range("F3").Value = -20.10157721253326
$ws.Range("G3").Value = -18.10157721253326
$ws.Range("H3").Value = -15.10157721253326
$ws.Range("I3").Value = -13.10157721253326
$ws.Range("J3").Value = -10.10157721253326
$ws.Range("K3").Value = -8.101577212533257
$ws.Range("L3").Value = -5.101577212533257
$ws.Range("M3").Value = -0.1015772125332575
$ws.Range("N3").Value = 4.898422787466743

# Row 5
$ws.Range("B5").Value = 73.552
$ws.Range("C5").Value = 70.62
$ws.Range("D5").Value = 67.887
$ws.Range("E5").Value = 62.974
$ws.Range("F5").Value = 120.589
$ws.Range("G5").Value = 87.306
$ws.Range("H5").Value = 122.819
$ws.Range("I5").Value = 89.946
$ws.Range("J5").Value = 56.812
$ws.Range("K5").Value = 56.086
$ws.Range("L5").Value = 55.035
$ws.Range("M5").Value = 53.377
$ws.Range("N5").Value = 51.83

# Row 6
$ws.Range("B6").Value = 36.085
$ws.Range("C6").Value = 34.646
$ws.Range("D6").Value = 33.305
$ws.Range("E6").Value = 30.895
$ws.Range("F6").Value = 59.161
$ws.Range("G6").Value = 42.832
$ws.Range("H6").Value = 60.255
$ws.Range("I6").Value = 44.128
$ws.Range("J6").Value = 27.872
$ws.Range("K6").Value = 27.516
$ws.Range("L6").Value = 27.0
$ws.Range("M6").Value = 26.187
$ws.Range("N6").Value = 25.428

# Row 7
$ws.Range("B7").Value = 2241.885
$ws.Range("C7").Value = 2152.517
$ws.Range("D7").Value = 2069.215
$ws.Range("E7").Value = 1919.465
$ws.Range("F7").Value = 1345.878
$ws.Range("G7").Value = 1545.029
$ws.Range("H7").Value = 1300.333
$ws.Range("I7").Value = 1474.329
$ws.Range("J7").Value = 1731.645
$ws.Range("K7").Value = 1709.517
$ws.Range("L7").Value = 1677.482
$ws.Range("M7").Value = 1626.946
$ws.Range("N7").Value = 1579.793

# Row 8
$ws.Range("B8").Value = 30.013
$ws.Range("C8").Value = 28.816
$ws.Range("D8").Value = 27.701
$ws.Range("E8").Value = 25.697
$ws.Range("F8").Value = 49.207
$ws.Range("G8").Value = 35.625
$ws.Range("H8").Value = 50.116
$ws.Range("I8").Value = 36.703
$ws.Range("J8").Value = 23.182
$ws.Range("K8").Value = 22.886
$ws.Range("L8").Value = 22.457
$ws.Range("M8").Value = 21.781
$ws.Range("N8").Value = 21.149

# Row 9
$ws.Range("B9").Value = 4.275
$ws.Range("C9").Value = 4.452
$ws.Range("D9").Value = 4.632
$ws.Range("E9").Value = 4.993
$ws.Range("F9").Value = 7.124
$ws.Range("G9").Value = 6.205
$ws.Range("H9").Value = 7.374
$ws.Range("I9").Value = 6.503
$ws.Range("J9").Value = 5.535
$ws.Range("K9").Value = 5.607
$ws.Range("L9").Value = 5.714
$ws.Range("M9").Value = 5.892
$ws.Range("N9").Value = 6.068

# Row 10
$ws.Range("B10").Value = 5.902
$ws.Range("C10").Value = 6.148
$ws.Range("D10").Value = 6.395
$ws.Range("E10").Value = 6.895
$ws.Range("F10").Value = 9.839
$ws.Range("G10").Value = 8.569
$ws.Range("H10").Value = 10.185
$ws.Range("I10").Value = 8.98
$ws.Range("J10").Value = 7.644
$ws.Range("K10").Value = 7.743
$ws.Range("L10").Value = 7.891
$ws.Range("M10").Value = 8.137
$ws.Range("N10").Value = 8.38

# Row 11
$ws.Range("B11").Value = 3.78
$ws.Range("C11").Value = 3.937
$ws.Range("D11").Value = 4.095
$ws.Range("E11").Value = 4.415
$ws.Range("F11").Value = 6.299
$ws.Range("G11").Value = 5.486
$ws.Range("H11").Value = 6.52
$ws.Range("I11").Value = 5.749
$ws.Range("J11").Value = 4.894
$ws.Range("K11").Value = 4.958
$ws.Range("L11").Value = 5.052
$ws.Range("M11").Value = 5.209
$ws.Range("N11").Value = 5.365

# Row 12
$ws.Range("B12").Value = 5.02
$ws.Range("C12").Value = 5.229
$ws.Range("D12").Value = 5.439
$ws.Range("E12").Value = 5.864
$ws.Range("F12").Value = 8.368
$ws.Range("G12").Value = 7.287
$ws.Range("H12").Value = 8.661
$ws.Range("I12").Value = 7.637
$ws.Range("J12").Value = 6.501
$ws.Range("K12").Value = 6.585
$ws.Range("L12").Value = 6.711
$ws.Range("M12").Value = 6.92
$ws.Range("N12").Value = 7.127

# Row 13
$ws.Range("B13").Value = 1.099
$ws.Range("C13").Value = 1.145
$ws.Range("D13").Value = 1.191
$ws.Range("E13").Value = 1.284
$ws.Range("F13").Value = 1.831
$ws.Range("G13").Value = 1.595
$ws.Range("H13").Value = 1.896
$ws.Range("I13").Value = 1.672
$ws.Range("J13").Value = 1.423
$ws.Range("K13").Value = 1.442
$ws.Range("L13").Value = 1.469
$ws.Range("M13").Value = 1.515
$ws.Range("N13").Value = 1.56
